$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A and B both end up at OOXML width 15.42578125.
# The runtime quantizes ColumnWidth to steps of 1/6, with stored = ColumnWidth + 5/6 (rounded).
# 14.666666666666666 is the closest achievable setting (-> stored width 15.5, nearest reachable value).
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666

# Update cell values A1:B32
$ws.Cells.Item(1, 1).Value = -0.10955972832100969
$ws.Cells.Item(1, 2).Value = 0.10916505688611267
$ws.Cells.Item(2, 1).Value = -0.088225278769406223
$ws.Cells.Item(2, 2).Value = 0.086915849620994301
$ws.Cells.Item(3, 1).Value = -0.037200998095794802
$ws.Cells.Item(3, 2).Value = 0.036844836257838409
$ws.Cells.Item(4, 1).Value = -0.028844836334137369
$ws.Cells.Item(4, 2).Value = 0.028529223039214457
$ws.Cells.Item(5, 1).Value = -0.025529223074140184
$ws.Cells.Item(5, 2).Value = 0.024466481168246013
$ws.Cells.Item(6, 1).Value = -0.0060342143235647683
$ws.Cells.Item(6, 2).Value = 0.0059369259925965423
$ws.Cells.Item(7, 1).Value = 0.0040630739124689619
$ws.Cells.Item(7, 2).Value = -0.0040754381718173249
$ws.Cells.Item(8, 1).Value = 0.014075438077646218
$ws.Cells.Item(8, 2).Value = -0.014088347424723757
$ws.Cells.Item(9, 1).Value = 0.016088347401566061
$ws.Cells.Item(9, 2).Value = -0.016101071209629225
$ws.Cells.Item(10, 1).Value = 0.018101071189761342
$ws.Cells.Item(10, 2).Value = -0.018100833733431188
$ws.Cells.Item(11, 1).Value = 0.021100833705235189
$ws.Cells.Item(11, 2).Value = -0.021104004929320297
$ws.Cells.Item(12, 1).Value = 0.024604004897704979
$ws.Cells.Item(12, 2).Value = -0.024665763920093564
$ws.Cells.Item(13, 1).Value = 0.0055942300412343826
$ws.Cells.Item(13, 2).Value = -0.0057562681654657055
$ws.Cells.Item(14, 1).Value = 0.013756268103388258
$ws.Cells.Item(14, 2).Value = -0.013837053537301358
$ws.Cells.Item(15, 1).Value = 0.014837053536194134
$ws.Cells.Item(15, 2).Value = -0.014903261835722326
$ws.Cells.Item(16, 1).Value = -0.0060334742906524852
$ws.Cells.Item(16, 2).Value = 0.0060031928197750162
$ws.Cells.Item(17, 1).Value = -0.0040031928282102669
$ws.Cells.Item(17, 2).Value = 0.003999999974310775
$ws.Cells.Item(18, 1).Value = -0.014939778168638895
$ws.Cells.Item(18, 2).Value = 0.014896286081139465
$ws.Cells.Item(19, 1).Value = -0.012091192285973662
$ws.Cells.Item(19, 2).Value = 0.012016207153954017
$ws.Cells.Item(20, 1).Value = -0.0080162071918756794
$ws.Cells.Item(20, 2).Value = 0.0080056336267855244
$ws.Cells.Item(21, 1).Value = -0.004005633665157049
$ws.Cells.Item(21, 2).Value = 0.0039999999612696513
$ws.Cells.Item(22, 1).Value = -0.045714851566132531
$ws.Cells.Item(22, 2).Value = 0.04550096933464598
$ws.Cells.Item(23, 1).Value = -0.040500969386868313
$ws.Cells.Item(23, 2).Value = 0.040099313992102203
$ws.Cells.Item(24, 1).Value = -0.020099314177199012
$ws.Cells.Item(24, 2).Value = 0.019999999812169378
$ws.Cells.Item(25, 1).Value = -0.016432266872422829
$ws.Cells.Item(25, 2).Value = 0.016372556826244278
$ws.Cells.Item(26, 1).Value = -0.01387255685847677
$ws.Cells.Item(26, 2).Value = 0.013798193821294547
$ws.Cells.Item(27, 1).Value = -0.011298193853911176
$ws.Cells.Item(27, 2).Value = 0.010872394386002071
$ws.Cells.Item(28, 1).Value = -0.0088723944157598211
$ws.Cells.Item(28, 2).Value = 0.0085996526416050045
$ws.Cells.Item(29, 1).Value = -0.0015996527143977701
$ws.Cells.Item(29, 2).Value = 0.0015304042510555504
$ws.Cells.Item(30, 1).Value = 0.058469595225975279
$ws.Cells.Item(30, 2).Value = -0.058816755775485063
$ws.Cells.Item(31, 1).Value = 0.065816755713235864
$ws.Cells.Item(31, 2).Value = -0.065940756056054539
$ws.Cells.Item(32, 1).Value = -0.0040010017351157501
$ws.Cells.Item(32, 2).Value = 0.0039999999665134567
